$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing raw-score inputs for Shafqat Nur (row 22)
$ws.Range("C22").Value = 14
$ws.Range("D22").Value = 2
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = 2
$ws.Range("K22").Value = 14
$ws.Range("L22").Value = 2
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("S22").Value = 9
$ws.Range("T22").Value = 6

$excel.Calculate()

# Update the view: scroll so column F is the left-most visible column,
# and select C22:T22 as the active selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("C22:T22").Select()
